# Adds a new column S (year 2022) to the "11.5.1 Number of deaths
# attributed to disasters" worksheet, mirroring the formatting already
# used by column R (year 2021) for rows 3-34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 2022 values (or "-" for "no data") for rows 5-34; row 4 is the
# 2022 header/year and row 3 is the blank separator row directly above it.
$values = @{
    4  = 2022
    5  = 135
    6  = 99
    7  = 36
    8  = 97
    9  = 80
    10 = 17
    11 = 17
    12 = 11
    13 = 6
    14 = 5
    15 = 3
    16 = 2
    17 = "-"
    18 = "-"
    19 = "-"
    20 = 6
    21 = 1
    22 = 5
    23 = "-"
    24 = "-"
    25 = "-"
    26 = 10
    27 = 4
    28 = 6
    29 = "-"
    30 = "-"
    31 = "-"
    32 = "-"
    33 = "-"
    34 = "-"
}

# Column S should look exactly like column R for every one of these
# rows, so copy the formatting across first ...
for ($r = 3; $r -le 34; $r++) {
    $src = $ws.Range("R" + $r)
    $dst = $ws.Range("S" + $r)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

# ... then fill in the actual 2022 figures (row 3 stays empty, matching
# the blank separator cells in column R).
foreach ($row in $values.Keys) {
    $ws.Range("S" + $row).Value = $values[$row]
}

$excel.CutCopyMode = $false

# Match the author's new selection/active cell.
[void]$ws.Range("S3").Select()
